$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

# --- Update comment texts ---

# O1: single_file_export_format comment, was "TODO"
$ws.Range("O1").Comment.Text("The format in which each single imaging file will be exported. (Example: DICOM, tiff, avi, etc.)")

# AF1: voltage_unit comment, "Example: volt" -> "Example: V"
$ws.Range("AF1").Comment.Text("The unit of the voltage used to acquire microCT images. Example: V")

# AH1: current_unit comment, "Example: amp" -> "Example: mA"
$ws.Range("AH1").Comment.Text("The unit of the current used to acquire microCT images. Example: mA")

# AL1: total_sections_analyzed comment, add "or OCT"
$ws.Range("AL1").Comment.Text("The number of sections used for analyzing microCT or OCT images")

# --- Update list sheet values (shared strings) ---

# voltage_unit list: "kVa" -> "kV"
$wsVoltage = $wb.Worksheets.Item("voltage_unit list")
$wsVoltage.Range("A2").Value = "kV"

# current_unit list: "mAmp" -> "mA", "microAmp" -> "microA", "Amp" -> "A"
$wsCurrent = $wb.Worksheets.Item("current_unit list")
$wsCurrent.Range("A1").Value = "mA"
$wsCurrent.Range("A2").Value = "microA"
$wsCurrent.Range("A3").Value = "A"

# --- Update data validation error messages on main sheet to match new list values ---

$dvVoltage = $ws.Range("AF2:AF1048576").Validation
$dvVoltage.ErrorMessage = "Value must be one of: V / kV."

$dvCurrent = $ws.Range("AH2:AH1048576").Validation
$dvCurrent.ErrorMessage = "Value must be one of: mA / microA / A."
